# Apply reordering of the goods-type labels in column A (shared-string table
# reorder) while leaving the Counts numbers in column B untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "деревенский товар"
$ws.Range("A18").Value = "серебреный товар"
$ws.Range("A21").Value = "железный товар"
$ws.Range("A22").Value = "мясо"
$ws.Range("A24").Value = "щепетильный товар"
$ws.Range("A25").Value = "пушной товар"
$ws.Range("A27").Value = "нужный товар"
$ws.Range("A28").Value = "питейный припасы"
$ws.Range("A29").Value = "внутренний товар"
$ws.Range("A30").Value = "медный товар"
$ws.Range("A32").Value = "оловянный товар"
$ws.Range("A33").Value = "привозный товар"
$ws.Range("A34").Value = "суровский товар"
$ws.Range("A36").Value = "заморский товар"
$ws.Range("A37").Value = "галантерейный товар"
$ws.Range("A39").Value = "харчевой припасы"
$ws.Range("A40").Value = "надлежащий товар"
$ws.Range("A41").Value = "меховой товар"
$ws.Range("A42").Value = "рукодельный товар"
$ws.Range("A43").Value = "домовый товар"
